$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that can look numeric (e.g. "211.81").
# Force the cell to stay a text value (matching the source inlineStr cells)
# and avoid picking up a Text-format style by resetting Style afterward.
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" '26.313.88'
$ws.Range("E2").Value = '  +0.50%  '
Set-TextValue "D3" '1.595.90'
$ws.Range("E3").Value = '  +0.29%  '
Set-TextValue "D5" '211.81'
$ws.Range("E5").Value = '  -0.14%  '
$ws.Range("E6").Value = '  -0.10%  '
$ws.Range("E8").Value = '  -0.15%  '
$ws.Range("E9").Value = '  -0.03%  '
Set-TextValue "D10" '19.07'
$ws.Range("E10").Value = '  +0.48%  '
$ws.Range("E11").Value = '  +1.11%  '
Set-TextValue "D12" '1.820.07'
$ws.Range("E12").Value = '  +0.28%  '
Set-TextValue "D13" '1.589.61'
$ws.Range("E13").Value = '  -0.07%  '
$ws.Range("E14").Value = '  -0.68%  '
$ws.Range("E15").Value = '  -1.04%  '
Set-TextValue "D16" '63.45'
$ws.Range("E16").Value = '  -0.30%  '
Set-TextValue "D17" '26.307.60'
$ws.Range("E17").Value = '  +0.47%  '
Set-TextValue "D18" '229.76'
$ws.Range("E18").Value = '  +7.47%  '
Set-TextValue "D19" '7.65'
$ws.Range("E19").Value = '  +4.37%  '
Set-TextValue "D20" '0.0₃0722'
$ws.Range("E20").Value = '  -0.44%  '
$ws.Range("E21").Value = '  +0.06%  '
$ws.Range("E22").Value = '  -0.22%  '
$ws.Range("E23").Value = '  +2.74%  '
$ws.Range("E24").Value = '  -1.24%  '
Set-TextValue "D25" '146.44'
$ws.Range("E25").Value = '  +1.22%  '
$ws.Range("E26").Value = '  +0.01%  '
Set-TextValue "D27" '6.97'
$ws.Range("E27").Value = '  +0.14%  '
$ws.Range("E28").Value = '  +0.30%  '
Set-TextValue "D29" '15.35'
$ws.Range("E29").Value = '  +1.79%  '
$ws.Range("E30").Value = '  +0.16%  '
$ws.Range("E31").Value = '  -0.07%  '
Set-TextValue "D32" '1.496.10'
$ws.Range("E32").Value = '  +5.29%  '
Set-TextValue "D33" '3.21'
$ws.Range("E33").Value = '  +1.29%  '
$ws.Range("E34").Value = '  -0.87%  '
$ws.Range("E35").Value = '  -0.26%  '
$ws.Range("E36").Value = '  +0.56%  '
Set-TextValue "D37" '0.569'
$ws.Range("E37").Value = '  -3.02%  '
$ws.Range("E38").Value = '  -0.77%  '
$ws.Range("E39").Value = '  -0.54%  '
$ws.Range("E40").Value = '  -1.89%  '
$ws.Range("E41").Value = '  +0.07%  '
$ws.Range("B42").Value = 'WEMIXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue "D42" '0.941'
$ws.Range("E42").Value = '  -4.75%  '
$ws.Range("B43").Value = 'MXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue "D43" '2.16'
$ws.Range("E43").Value = '  +1.87%  '
Set-TextValue "D44" '1.733.10'
$ws.Range("E45").Value = '  -0.86%  '
$ws.Range("E46").Value = '  -0.55%  '
$ws.Range("E47").Value = '  +1.68%  '
$ws.Range("E48").Value = '  -0.21%  '
Set-TextValue "D49" '0.0501'
$ws.Range("E49").Value = '  -0.15%  '
$ws.Range("E50").Value = '  +0.11%  '
$ws.Range("E51").Value = '  +0.09%  '
